$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like columns that could be misinterpreted as numbers/dates by Excel
# auto-conversion keep their literal text representation.
$ws.Range("Y2:AB24").NumberFormat = "@"
$ws.Range("I2:I24").NumberFormat = "@"

$rows = @(
  @{ Row=2; "A"=111896635; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575037.2974304935; "R"=6703389.027347369; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=3; "A"=111896690; "B"=90687; "C"="Ovaliderad"; "D"="LC"; "E"=5964; "F"="Fjällig taggsvamp s.str."; "G"="Sarcodon imbricatus s.str."; "H"="(L.:Fr.) P.Karst."; "P"="Kratte masugn, Gstr"; "Q"=575060.2881161601; "R"=6703376.67477417; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=4; "A"=111896638; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575087.1320314853; "R"=6703393.020834555; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=5; "A"=111884093; "B"=98535; "C"="Ovaliderad"; "D"="LC"; "E"=222498; "F"="Blåsippa"; "G"="Hepatica nobilis"; "H"="Schreb."; "P"="Kopparåsen (Kopparåsen), Gstr"; "Q"=575065.9914513066; "R"=6703387.648325931; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Patric Engfeldt"; "AX"="Patric Engfeldt" },
  @{ Row=6; "A"=111896641; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575021.3626164712; "R"=6703370.933926445; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=7; "A"=111896637; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575088.0587098968; "R"=6703396.00058554; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=8; "A"=111896640; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575025.3556637274; "R"=6703369.042946251; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=9; "A"=111896633; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575100.4050603262; "R"=6703444.118284944; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=10; "A"=111896639; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575089.384229039; "R"=6703379.745088123; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=11; "A"=111896636; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575108.85141061; "R"=6703418.142308297; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=12; "A"=111884133; "B"=88899; "C"="Ovaliderad"; "D"="NT"; "E"=3286; "F"="Flattoppad klubbsvamp"; "G"="Clavariadelphus truncatus"; "H"="(Quél.) Donk"; "P"="Kalkberget (Kalkberget), Gstr"; "Q"=575059.034285416; "R"=6703389.477814267; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Patric Engfeldt"; "AX"="Patric Engfeldt" },
  @{ Row=13; "A"=111883983; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kalkberget (Kalkberget), Gstr"; "Q"=575058.3527020445; "R"=6703446.206921679; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Patric Engfeldt"; "AX"="Patric Engfeldt" },
  @{ Row=14; "A"=111896653; "B"=89183; "C"="Ovaliderad"; "D"="LC"; "E"=3215; "F"="Rödgul trumpetsvamp"; "G"="Craterellus lutescens"; "H"="(Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575075.050630242; "R"=6703403.625642136; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=15; "A"=111896643; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575038.7114136803; "R"=6703416.194821274; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=16; "A"=111896655; "B"=89183; "C"="Ovaliderad"; "D"="LC"; "E"=3215; "F"="Rödgul trumpetsvamp"; "G"="Craterellus lutescens"; "H"="(Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575104.6742508161; "R"=6703428.910891063; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=17; "A"=111896634; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575048.3395925189; "R"=6703452.413791304; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=18; "A"=111896642; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575014.1091647458; "R"=6703387.066676207; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=19; "A"=111884471; "B"=88899; "C"="Ovaliderad"; "D"="NT"; "E"=3286; "F"="Flattoppad klubbsvamp"; "G"="Clavariadelphus truncatus"; "H"="(Quél.) Donk"; "P"="Kalkberget (Kalkberget), Gstr"; "Q"=575020.8210917887; "R"=6703397.074168184; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Patric Engfeldt"; "AX"="Patric Engfeldt" },
  @{ Row=20; "A"=111896652; "B"=89183; "C"="Ovaliderad"; "D"="LC"; "E"=3215; "F"="Rödgul trumpetsvamp"; "G"="Craterellus lutescens"; "H"="(Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575066.556649723; "R"=6703455.751857814; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=21; "A"=111896654; "B"=89183; "C"="Ovaliderad"; "D"="LC"; "E"=3215; "F"="Rödgul trumpetsvamp"; "G"="Craterellus lutescens"; "H"="(Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575072.6962527435; "R"=6703421.833381963; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=22; "A"=111896644; "B"=90332; "C"="Ovaliderad"; "D"="LC"; "E"=4769; "F"="Svavelriska"; "G"="Lactarius scrobiculatus"; "H"="(Scop.:Fr.) Fr."; "P"="Kratte masugn, Gstr"; "Q"=575036.4083237475; "R"=6703431.936489306; "S"=25; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-04"; "Z"="00:00"; "AA"="2023-09-04"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=23; "A"=112204545; "B"=90687; "C"="Ovaliderad"; "D"="LC"; "E"=5964; "F"="Fjällig taggsvamp s.str."; "G"="Sarcodon imbricatus s.str."; "H"="(L.:Fr.) P.Karst."; "P"="Kratte Masugn, Gstr"; "Q"=575051.354848919; "R"=6703378.463325701; "S"=15; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-18"; "Z"="00:00"; "AA"="2023-09-18"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" },
  @{ Row=24; "A"=112204541; "B"=88966; "C"="Ovaliderad"; "D"="NT"; "E"=5754; "F"="Gultoppig fingersvamp"; "G"="Ramaria testaceoflava"; "H"="(Bres.) Corner"; "I"="7"; "J"="fruktkroppar"; "P"="Kratte Masugn, Gstr"; "Q"=575058.4344827051; "R"=6703442.261289957; "S"=15; "T"="Gävleborg"; "U"="Hofors"; "V"="Gästrikland"; "W"="Torsåker"; "Y"="2023-09-18"; "Z"="00:00"; "AA"="2023-09-18"; "AB"="00:00"; "AD"=$false; "AE"=$false; "AG"=$false; "AW"="Philipp Weiss"; "AX"="Philipp Weiss" }
)

foreach ($r in $rows) {
    $rownum = $r.Row
    foreach ($key in $r.Keys) {
        if ($key -eq "Row") { continue }
        $addr = "$key$rownum"
        $ws.Range($addr).Value = $r[$key]
    }
}

Write-Host "Done."